$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C8"   = "[name=`"Folinic`"]  Just 'hold them accountable?'`n"
    "C10"  = "[name=`"Folinic`"]  Do you have anything at all to tell 'a couple of tourists' then? 'Schultz?'`n"
    "C13"  = "[name=`"Folinic`"]  Usually, we go by 'Operator.'`n"
    "C37"  = "[name=`"Folinic`"]  The Infected rebels are acting really weird. And they keep talking about a 'fire' and a 'doctor,' right?`n"
    "C56"  = "[name=`"Folinic`"]  Where we come from, most people prefer not to use their real names at work. Just like 'Atro.'`n"
    "C78"  = "[name=`"Folinic`"]  The local Infected seem to know 'Dr. Atro.' You’re stalling, and I’d prefer you not.`n"
    "C84"  = "[name=`"Severin`"]  We received the Infected doctor who called herself 'Atro.' And we greatly appreciate her help caring for the many Infected residents of Wolumonde during this difficult time.`n"
    "C85"  = "[name=`"Folinic`"]  And you show your appreciation by telling us 'she disappeared?'`n"
    "C101" = "[name=`"Severin`"]  If I’m not going to talk lightly, all I can give you is 'I’m very sorry.'`n"
    "C109" = "[name=`"Folinic`"]  The term is 'Sarkaz,' please. I hope the meaning of the word is not lost on you. Do you have any other information about the armed Infected?`n"
    "C111" = "[name=`"Folinic`"]  Oh? So why didn't they immediately show up to plunder the town? Why are they just 'loitering' out there?`n"
    "C119" = "[name=`"Folinic`"]  For example... have you ever heard of the 'Reunion Movement?'`n"
    "C127" = "[name=`"Folinic`"]  We may have to do some investigating around town. We need to find out what the deal with this 'disappearance' is...`n"
    "C153" = "[name=`"Severin`"]  They’re 'operators.' I only listen to 'doctors.'`n"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
